# "Se eleva el Trabajo POM finalizado"
# Rename the sheet, fix up several test-data cells (CP008 now targets the
# "Vinilo" product instead of "Arte"; a couple of login-test values were
# swapped/corrected), and turn the e-mail in B12 into a mailto hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Hoja1" -> "DatosCP"
$ws.Name = "DatosCP"

# CP008 test case: was "producto Arte", now "producto Vinilo"
$ws.Range("A9").Value = "CP008_productoViniloMenu"
$ws.Range("B9").Value = "Vinilo Puerta Dragón"

# CP010 row: expected password value corrected
$ws.Range("C11").Value = "Prueba"

# CP011 row (B12): turn the plain e-mail address into a mailto hyperlink.
# Excel applies the built-in "Hyperlink" cell style automatically.
$ws.Hyperlinks.Add($ws.Range("B12"), "mailto:maildeprueba29@noexite.com", [Type]::Missing, [Type]::Missing, "maildeprueba29@")
# The displayed text then gets trimmed down to the user part only.
$ws.Range("B12").Value = "maildeprueba29"

# CP012 row: email/password values updated
$ws.Range("B13").Value = "maildeprueba29@noexite.com"
$ws.Range("C13").Value = "Prueba123"

# Final selection left on C13
$ws.Range("C13").Select() | Out-Null
